$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("A1:N1").Style = "Normal"
$ws.Range("H6").Value = 127286.375
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H19").Value = 1433.0476
$ws.Range("I19").Value = 730.7
$ws.Range("J19").Value = 2071.5454
$ws.Range("K19").Value = 730.7
$ws.Range("L19").Value = 2071.5454
$ws.Range("M19").Value = -555.7
$ws.Range("N19").Value = -2421.5454
$ws.Range("H33").Value = 430.33334
$ws.Range("I33").Value = 144.5
$ws.Range("J33").Value = 1002
$ws.Range("K33").Value = 144.5
$ws.Range("L33").Value = 1002
$ws.Range("M33").Value = 84.5
$ws.Range("N33").Value = -1460
$ws.Range("H38").Value = 308.44446
$ws.Range("I38").Value = 308.44446
$ws.Range("K38").Value = 925.33338
$ws.Range("M38").Value = -553.33338
$ws.Range("H51").Value = 8409.622
$ws.Range("I51").Value = 7077.6
$ws.Range("J51").Value = 8617.75
$ws.Range("K51").Value = 7077.6
$ws.Range("L51").Value = 8617.75
$ws.Range("M51").Value = -6593.6
$ws.Range("N51").Value = -9585.75
$ws.Range("H111").Value = 29999
$ws.Range("I111").Value = 29999
$ws.Range("K111").Value = 89997
$ws.Range("M111").Value = -86930
$ws.Range("H113").Value = 2949.8147
$ws.Range("I113").Value = 3397.3157
$ws.Range("K113").Value = 3397.3157
$ws.Range("M113").Value = -143.3157000000001
$ws.Range("H116").Value = 5016.6875
$ws.Range("I116").Value = 4899.0835
$ws.Range("K116").Value = 4899.0835
$ws.Range("M116").Value = -1457.0835
$ws.Range("H121").Value = 3811.3333
$ws.Range("J121").Value = 3811.3333
$ws.Range("L121").Value = 11433.9999
$ws.Range("N121").Value = -14927.9999
$ws.Range("H125").Value = 981.8
$ws.Range("I125").Value = 981.8
$ws.Range("K125").Value = 8836.199999999999
$ws.Range("M125").Value = -6376.199999999999
$ws.Range("H132").Value = 2030
$ws.Range("I132").Value = 1724.6279
$ws.Range("J132").Value = 8595.5
$ws.Range("K132").Value = 5173.8837
$ws.Range("L132").Value = 25786.5
$ws.Range("M132").Value = -2643.8837
$ws.Range("N132").Value = -30846.5
$ws.Range("H138").Value = 4053.2942
$ws.Range("I138").Value = 3498.5
$ws.Range("J138").Value = 4172.1787
$ws.Range("K138").Value = 10495.5
$ws.Range("L138").Value = 12516.5361
$ws.Range("M138").Value = -5355.5
$ws.Range("N138").Value = -22796.5361
$ws.Range("H141").Value = 970.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("A1:N1").Style = "Normal"
$ws.Range("H2").Value = 2479.5
$ws.Range("I2").Value = 2561.875
$ws.Range("K2").Value = 2561.875
$ws.Range("M2").Value = -2448.875
$ws.Range("H31").Value = 4333
$ws.Range("I31").Value = 4333
$ws.Range("K31").Value = 4333
$ws.Range("M31").Value = -4039
$ws.Range("H32").Value = 2512.2678
$ws.Range("I32").Value = 2573.3396
$ws.Range("J32").Value = 1433.3334
$ws.Range("K32").Value = 2573.3396
$ws.Range("L32").Value = 1433.3334
$ws.Range("M32").Value = -2286.3396
$ws.Range("N32").Value = -2007.3334
$ws.Range("H63").Value = 2051.4546
$ws.Range("I63").Value = 2308
$ws.Range("K63").Value = 2308
$ws.Range("M63").Value = -1622
$ws.Range("H66").Value = 2051.4546
$ws.Range("I66").Value = 2308
$ws.Range("K66").Value = 11540
$ws.Range("M66").Value = -8108
$ws.Range("H74").Value = 111829.89
$ws.Range("I74").Value = 133216.06
$ws.Range("K74").Value = 133216.06
$ws.Range("M74").Value = -132342.06
$ws.Range("H77").Value = 111829.89
$ws.Range("I77").Value = 133216.06
$ws.Range("K77").Value = 666080.3
$ws.Range("M77").Value = -661712.3
$ws.Range("H102").Value = 1765.3334
$ws.Range("I102").Value = 1696.4
$ws.Range("K102").Value = 1696.4
$ws.Range("M102").Value = -74.40000000000009
$ws.Range("H110").Value = 1800.6136
$ws.Range("I110").Value = 1724.5135
$ws.Range("K110").Value = 1724.5135
$ws.Range("M110").Value = 320.4865
$ws.Range("H116").Value = 2479.5
$ws.Range("I116").Value = 2561.875
$ws.Range("K116").Value = 2561.875
$ws.Range("M116").Value = -267.875
$ws.Range("H132").Value = 1913.129
$ws.Range("I132").Value = 2028.8431
$ws.Range("K132").Value = 6086.5293
$ws.Range("M132").Value = -3556.5293

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("A1:N1").Style = "Normal"
$ws.Range("H3").Value = 2479.5
$ws.Range("I3").Value = 2561.875
$ws.Range("K3").Value = 2561.875
$ws.Range("M3").Value = -2447.875
$ws.Range("H94").Value = 675.1786
$ws.Range("I94").Value = 502.78262
$ws.Range("K94").Value = 502.78262
$ws.Range("M94").Value = -51.78262000000001
$ws.Range("H99").Value = 5030.1577
$ws.Range("J99").Value = 4994.5713
$ws.Range("L99").Value = 4994.5713
$ws.Range("N99").Value = -7990.5713
$ws.Range("H107").Value = 1228.6666
$ws.Range("I107").Value = 675.2
$ws.Range("J107").Value = 3996
$ws.Range("K107").Value = 675.2
$ws.Range("L107").Value = 3996
$ws.Range("M107").Value = 1244.8
$ws.Range("N107").Value = -7836
$ws.Range("H134").Value = 3191.6765
$ws.Range("I134").Value = 3883.353
$ws.Range("K134").Value = 11650.059
$ws.Range("M134").Value = -9115.059000000001
$ws.Range("H141").Value = 64000
$ws.Range("J141").Value = 64000
$ws.Range("L141").Value = 64000
$ws.Range("N141").Value = -74360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("A1:N1").Style = "Normal"
$ws.Range("H16").Value = 1265.1111
$ws.Range("I16").Value = 721.4762
$ws.Range("J16").Value = 3167.8333
$ws.Range("K16").Value = 721.4762
$ws.Range("L16").Value = 3167.8333
$ws.Range("M16").Value = -434.4761999999999
$ws.Range("N16").Value = -3741.8333
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H31").Value = 295541.62
$ws.Range("I31").Value = 590011.8
$ws.Range("K31").Value = 590011.8
$ws.Range("M31").Value = -589716.8
$ws.Range("H34").Value = 295541.62
$ws.Range("I34").Value = 590011.8
$ws.Range("K34").Value = 590011.8
$ws.Range("M34").Value = -589809.8
$ws.Range("H44").Value = 10044
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H58").Value = 1644.1666
$ws.Range("I58").Value = 1630
$ws.Range("K58").Value = 1630
$ws.Range("M58").Value = -1427
$ws.Range("H98").Value = 42371
$ws.Range("J98").Value = 42371
$ws.Range("L98").Value = 42371
$ws.Range("N98").Value = -46863
$ws.Range("H113").Value = 1265.1111
$ws.Range("I113").Value = 721.4762
$ws.Range("J113").Value = 3167.8333
$ws.Range("K113").Value = 721.4762
$ws.Range("L113").Value = 3167.8333
$ws.Range("M113").Value = 1448.5238
$ws.Range("N113").Value = -7507.8333
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H122").Value = 1202.4615
$ws.Range("I122").Value = 1177.6666
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3532.9998
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1082.9998
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 5150.05
$ws.Range("I132").Value = 5003.0884
$ws.Range("K132").Value = 15009.2652
$ws.Range("M132").Value = -12479.2652
$ws.Range("H134").Value = 2547.8262
$ws.Range("I134").Value = 2280.05
$ws.Range("J134").Value = 4333
$ws.Range("K134").Value = 6840.150000000001
$ws.Range("L134").Value = 12999
$ws.Range("M134").Value = -4305.150000000001
$ws.Range("N134").Value = -18069
$ws.Range("H136").Value = 1644.1666
$ws.Range("I136").Value = 1630
$ws.Range("K136").Value = 4890
$ws.Range("M136").Value = -2340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("A1:N1").Style = "Normal"
$ws.Range("H2").Value = 219.86667
$ws.Range("I2").Value = 297.42856
$ws.Range("J2").Value = 152
$ws.Range("K2").Value = 1784.57136
$ws.Range("L2").Value = 912
$ws.Range("M2").Value = -1671.57136
$ws.Range("N2").Value = -1138
$ws.Range("H13").Value = 123.6
$ws.Range("I13").Value = 152
$ws.Range("J13").Value = 10
$ws.Range("K13").Value = 456
$ws.Range("L13").Value = 30
$ws.Range("M13").Value = -288
$ws.Range("N13").Value = -366
$ws.Range("H34").Value = 1310.0667
$ws.Range("I34").Value = 74.38461
$ws.Range("J34").Value = 2255
$ws.Range("K34").Value = 223.15383
$ws.Range("L34").Value = 6765
$ws.Range("M34").Value = -139.15383
$ws.Range("N34").Value = -6933
$ws.Range("H39").Value = 1008.8333
$ws.Range("I39").Value = 828.7273
$ws.Range("K39").Value = 2486.1819
$ws.Range("M39").Value = -2192.1819
$ws.Range("H55").Value = 8989.286
$ws.Range("J55").Value = 8989.286
$ws.Range("L55").Value = 26967.858
$ws.Range("N55").Value = -27321.858
$ws.Range("H107").Value = 939.6071
$ws.Range("I107").Value = 459.18182
$ws.Range("K107").Value = 1377.54546
$ws.Range("M107").Value = 542.45454
$ws.Range("H109").Value = 201105.44
$ws.Range("I109").Value = 201105.44
$ws.Range("K109").Value = 603316.3200000001
$ws.Range("M109").Value = -602276.3200000001
$ws.Range("H131").Value = 5557199
$ws.Range("I131").Value = 100001100
$ws.Range("K131").Value = 300003300
$ws.Range("M131").Value = -299998260
$ws.Range("H132").Value = 5448.0938
$ws.Range("I132").Value = 7178.1904
$ws.Range("J132").Value = 2145.182
$ws.Range("K132").Value = 64603.7136
$ws.Range("L132").Value = 19306.638
$ws.Range("M132").Value = -62073.7136
$ws.Range("N132").Value = -24366.638
$ws.Range("H141").Value = 5500
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("A1:N1").Style = "Normal"
$ws.Range("H36").Value = 90999
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 37495
$ws.Range("J52").Value = 37495
$ws.Range("L52").Value = 37495
$ws.Range("N52").Value = -38013
$ws.Range("H80").Value = 3101.6667
$ws.Range("I80").Value = 2652.5
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2652.5
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1654.5
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 3101.6667
$ws.Range("I83").Value = 2652.5
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 13262.5
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -8270.5
$ws.Range("N83").Value = -29984
$ws.Range("H97").Value = 1946.5
$ws.Range("I97").Value = 1830.6
$ws.Range("J97").Value = 2139.6667
$ws.Range("K97").Value = 1830.6
$ws.Range("L97").Value = 2139.6667
$ws.Range("M97").Value = -1334.6
$ws.Range("N97").Value = -3131.6667
$ws.Range("H102").Value = 2729.4517
$ws.Range("I102").Value = 2379.0715
$ws.Range("K102").Value = 2379.0715
$ws.Range("M102").Value = -757.0715
$ws.Range("H132").Value = 61841.65
$ws.Range("I132").Value = 86076.164
$ws.Range("J132").Value = 3678.8
$ws.Range("K132").Value = 258228.492
$ws.Range("L132").Value = 11036.4
$ws.Range("M132").Value = -255698.492
$ws.Range("N132").Value = -16096.4
$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -95070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("A1:N1").Style = "Normal"
$ws.Range("H7").Value = 6340.636
$ws.Range("I7").Value = 6340.636
$ws.Range("K7").Value = 6340.636
$ws.Range("M7").Value = -6228.636
$ws.Range("H13").Value = 1503
$ws.Range("I13").Value = 6
$ws.Range("K13").Value = 6
$ws.Range("M13").Value = 134
$ws.Range("H20").Value = 9965.272
$ws.Range("J20").Value = 9965.272
$ws.Range("L20").Value = 9965.272
$ws.Range("N20").Value = -10417.272
$ws.Range("H46").Value = 3355.9583
$ws.Range("I46").Value = 2993.625
$ws.Range("J46").Value = 4080.625
$ws.Range("K46").Value = 2993.625
$ws.Range("L46").Value = 4080.625
$ws.Range("M46").Value = -2805.625
$ws.Range("N46").Value = -4456.625
$ws.Range("H61").Value = 1374.75
$ws.Range("I61").Value = 1374.75
$ws.Range("K61").Value = 1374.75
$ws.Range("M61").Value = -1172.75
$ws.Range("H82").Value = 1713.52
$ws.Range("I82").Value = 1481.591
$ws.Range("J82").Value = 3414.3333
$ws.Range("K82").Value = 1481.591
$ws.Range("L82").Value = 3414.3333
$ws.Range("M82").Value = -1120.591
$ws.Range("N82").Value = -4136.3333
$ws.Range("H85").Value = 1713.52
$ws.Range("I85").Value = 1481.591
$ws.Range("J85").Value = 3414.3333
$ws.Range("K85").Value = 1481.591
$ws.Range("L85").Value = 3414.3333
$ws.Range("M85").Value = -233.5909999999999
$ws.Range("N85").Value = -5910.3333
$ws.Range("H93").Value = 3999.5
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 3999.5
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 3999.5
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -6495.5
$ws.Range("H100").Value = 16856.428
$ws.Range("J100").Value = 14999.5
$ws.Range("L100").Value = 14999.5
$ws.Range("N100").Value = -16081.5
$ws.Range("H113").Value = 1374.75
$ws.Range("I113").Value = 1374.75
$ws.Range("K113").Value = 1374.75
$ws.Range("M113").Value = 795.25
$ws.Range("H126").Value = 6340.636
$ws.Range("I126").Value = 6340.636
$ws.Range("K126").Value = 19021.908
$ws.Range("M126").Value = -16551.908
$ws.Range("H132").Value = 3883.6875
$ws.Range("I132").Value = 3905.5715
$ws.Range("J132").Value = 3730.5
$ws.Range("K132").Value = 11716.7145
$ws.Range("L132").Value = 11191.5
$ws.Range("M132").Value = -9186.7145
$ws.Range("N132").Value = -16251.5
$ws.Range("H136").Value = 3501.6667
$ws.Range("I136").Value = 3549.5715
$ws.Range("J136").Value = 3166.3333
$ws.Range("K136").Value = 10648.7145
$ws.Range("L136").Value = 9498.999899999999
$ws.Range("M136").Value = -8098.7145
$ws.Range("N136").Value = -14598.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("A1:N1").Style = "Normal"
$ws.Range("H51").Value = 9997.6
$ws.Range("I51").Value = 9997.6
$ws.Range("K51").Value = 9997.6
$ws.Range("M51").Value = -9487.6
$ws.Range("H113").Value = 672.34485
$ws.Range("I113").Value = 667.3333
$ws.Range("K113").Value = 2001.9999
$ws.Range("M113").Value = 168.0001
$ws.Range("H126").Value = 3694.7
$ws.Range("I126").Value = 3333.3333
$ws.Range("J126").Value = 4236.75
$ws.Range("K126").Value = 9999.999899999999
$ws.Range("L126").Value = 12710.25
$ws.Range("M126").Value = -7529.999899999999
$ws.Range("N126").Value = -17650.25
$ws.Range("H132").Value = 1831.0857
$ws.Range("I132").Value = 1856.0588
$ws.Range("K132").Value = 5568.1764
$ws.Range("M132").Value = -3038.1764
$ws.Range("H136").Value = 1002796.7
$ws.Range("I136").Value = 1113940.8
$ws.Range("K136").Value = 3341822.4
$ws.Range("M136").Value = -3339272.4

